$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 102.4929963333333
$ws.Range("H2").Value = 307.478989
$ws.Range("I2").Value = 0.2065071987599813
$ws.Range("J2").Value = 0.2065071987599814
$ws.Range("M2").Value = 68.46613766666667
$ws.Range("N2").Value = 205.398413
$ws.Range("O2").Value = 0.4719163120948675
$ws.Range("P2").Value = 0.4719163120948675
$ws.Range("Q2").Value = 7017.299596827163
$ws.Range("R2").Value = 63155.69637144446
$ws.Range("S2").Value = 0.09745411565985218
$ws.Range("T2").Value = 0.09745411565985218

# Row 3
$ws.Range("G3").Value = 102.4929963333333
$ws.Range("H3").Value = 307.478989
$ws.Range("I3").Value = 0.2065071987599813
$ws.Range("J3").Value = 0.2065071987599814
$ws.Range("M3").Value = 9.278736333333333
$ws.Range("N3").Value = 27.836209
$ws.Range("O3").Value = 0.06395551407683932
$ws.Range("P3").Value = 0.06395551407683933
$ws.Range("Q3").Value = 951.0054889903001
$ws.Range("R3").Value = 8559.049400912701
$ws.Range("S3").Value = 0.01320727405726264
$ws.Range("T3").Value = 0.01320727405726265

# Row 4
$ws.Range("G4").Value = 102.4929963333333
$ws.Range("H4").Value = 307.478989
$ws.Range("I4").Value = 0.2065071987599813
$ws.Range("J4").Value = 0.2065071987599814
$ws.Range("M4").Value = 67.336226
$ws.Range("N4").Value = 202.008678
$ws.Range("O4").Value = 0.4641281738282933
$ws.Range("P4").Value = 0.4641281738282933
$ws.Range("Q4").Value = 6901.491564518505
$ws.Range("R4").Value = 62113.42408066655
$ws.Range("S4").Value = 0.09584580904286653
$ws.Range("T4").Value = 0.09584580904286655

# Row 5
$ws.Range("I5").Value = 0.581825957350084
$ws.Range("J5").Value = 0.5818259573500841
$ws.Range("M5").Value = 68.46613766666667
$ws.Range("N5").Value = 205.398413
$ws.Range("O5").Value = 0.4719163120948675
$ws.Range("P5").Value = 0.4719163120948675
$ws.Range("Q5").Value = 19770.96721302062
$ws.Range("R5").Value = 177938.7049171855
$ws.Range("S5").Value = 0.2745731600737173
$ws.Range("T5").Value = 0.2745731600737173

# Row 6
$ws.Range("I6").Value = 0.581825957350084
$ws.Range("J6").Value = 0.5818259573500841
$ws.Range("M6").Value = 9.278736333333333
$ws.Range("N6").Value = 27.836209
$ws.Range("O6").Value = 0.06395551407683932
$ws.Range("P6").Value = 0.06395551407683933
$ws.Range("Q6").Value = 2679.420777578205
$ws.Range("R6").Value = 24114.78699820385
$ws.Range("S6").Value = 0.03721097820557381
$ws.Range("T6").Value = 0.03721097820557383

# Row 7
$ws.Range("I7").Value = 0.581825957350084
$ws.Range("J7").Value = 0.5818259573500841
$ws.Range("M7").Value = 67.336226
$ws.Range("N7").Value = 202.008678
$ws.Range("O7").Value = 0.4641281738282933
$ws.Range("P7").Value = 0.4641281738282933
$ws.Range("Q7").Value = 19444.68261049144
$ws.Range("R7").Value = 175002.143494423
$ws.Range("S7").Value = 0.2700418190707929
$ws.Range("T7").Value = 0.270041819070793

# Row 8
$ws.Range("G8").Value = 105.053815
$ws.Range("H8").Value = 315.161445
$ws.Range("I8").Value = 0.2116668438899346
$ws.Range("J8").Value = 0.2116668438899346
$ws.Range("M8").Value = 68.46613766666667
$ws.Range("N8").Value = 205.398413
$ws.Range("O8").Value = 0.4719163120948675
$ws.Range("P8").Value = 0.4719163120948675
$ws.Range("Q8").Value = 7192.628960198532
$ws.Range("R8").Value = 64733.66064178679
$ws.Range("S8").Value = 0.09988903636129798
$ws.Range("T8").Value = 0.09988903636129799

# Row 9
$ws.Range("G9").Value = 105.053815
$ws.Range("H9").Value = 315.161445
$ws.Range("I9").Value = 0.2116668438899346
$ws.Range("J9").Value = 0.2116668438899346
$ws.Range("M9").Value = 9.278736333333333
$ws.Range("N9").Value = 27.836209
$ws.Range("O9").Value = 0.06395551407683932
$ws.Range("P9").Value = 0.06395551407683933
$ws.Range("Q9").Value = 974.7666501957783
$ws.Range("R9").Value = 8772.899851762006
$ws.Range("S9").Value = 0.01353726181400286
$ws.Range("T9").Value = 0.01353726181400287

# Row 10
$ws.Range("G10").Value = 105.053815
$ws.Range("H10").Value = 315.161445
$ws.Range("I10").Value = 0.2116668438899346
$ws.Range("J10").Value = 0.2116668438899346
$ws.Range("M10").Value = 67.336226
$ws.Range("N10").Value = 202.008678
$ws.Range("O10").Value = 0.4641281738282933
$ws.Range("P10").Value = 0.4641281738282933
$ws.Range("Q10").Value = 7073.927429002189
$ws.Range("R10").Value = 63665.34686101972
$ws.Range("S10").Value = 0.0982405457146338
$ws.Range("T10").Value = 0.09824054571463381
